# "Generate Report for Archive"
#
# The localization status for the e2e test assets moved on from the
# handoff stage, so every cell that still shows the old "Ready for
# handoff" status now reads "In Translation". That status column is
# narrower than the old text needed, so its column(s) are re-sized to
# fit the shorter label on every sheet that carries a status column.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = $newStatus
$overview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5

# --- Per-language detail sheets: column C holds the status ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2:C4").Value = $newStatus
    $ws.Range("C1").EntireColumn.ColumnWidth = 12.5
}
